# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 16:50"

# Update "Recuperados" (column D) values for several Canary Islands provinces
$ws.Range("D32").Value = 1056
$ws.Range("D43").Value = 1056
$ws.Range("D54").Value = 1056
$ws.Range("D56").Value = 1056
$ws.Range("D57").Value = 1056
$ws.Range("D61").Value = 1056
$ws.Range("D63").Value = 1056
